$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B2 from 33 to 66
$ws.Range("B2").Value = 66

# Add new row 5 with data
$ws.Range("A5").Value = "Cereza"
$ws.Range("B5").Value = 33
$ws.Range("C5").Value = 1.1
